# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the file
# "819d26cd-ee36-4459-b2d5-055c55a18cc3" is now ready for handoff
# (machine translation priority), across the Overview, zh-cn and de-de
# sheets, and widens the affected "Status"/language columns to fit the
# new text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (row 3 = 819d26cd-ee36-4459-b2d5-055c55a18cc3.md) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-01 00:16:10"
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3

# ---- zh-cn sheet (row 3 = 819d26cd-ee36-4459-b2d5-055c55a18cc3) ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-01 00:15:59"
$zhcn.Columns.Item(3).ColumnWidth = 16.3

# ---- de-de sheet (row 3 = 819d26cd-ee36-4459-b2d5-055c55a18cc3) ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-01 00:16:10"
$dede.Columns.Item(3).ColumnWidth = 16.3
